# Changing waits in suite B
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update the Runmode column (C) values from "N" to "Y" for rows 3-7
$ws.Range("C3:C7").Value = "Y"

# Update selection to C2
$ws.Range("C2").Select()
